$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells: Wins, Losses, Ties
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy formatting from the existing header cell (AC1) so the new headers
# match the bold/centered/bordered look of the rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the team record (Wins/Losses/Ties) for every data row (2-48)
for ($r = 2; $r -le 48; $r++) {
    $ws.Cells.Item($r, 30).Value = 80   # AD - Wins
    $ws.Cells.Item($r, 31).Value = 82   # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF - Ties
}
